$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Range("H6").Value = 2067361.1
$ws.Range("I6").Value = 2526674.8
$ws.Range("J6").Value = 450
$ws.Range("K6").Value = 7580024.399999999
$ws.Range("L6").Value = 1350
$ws.Range("M6").Value = -7579912.399999999
$ws.Range("N6").Value = -1574
$ws.Range("H8").Value = 300
$ws.Range("I8").Value = 300
$ws.Range("K8").Value = 900
$ws.Range("M8").Value = -761
$ws.Range("H17").Value = 670098.7
$ws.Range("J17").Value = 670098.7
$ws.Range("L17").Value = 2010296.1
$ws.Range("N17").Value = -2010632.1
$ws.Range("H34").Value = 4315
$ws.Range("I34").Value = 472.5
$ws.Range("J34").Value = 12000
$ws.Range("K34").Value = 472.5
$ws.Range("L34").Value = 12000
$ws.Range("M34").Value = -269.5
$ws.Range("N34").Value = -12406
$ws.Range("H36").Value = 4315
$ws.Range("I36").Value = 472.5
$ws.Range("J36").Value = 12000
$ws.Range("K36").Value = 472.5
$ws.Range("L36").Value = 12000
$ws.Range("M36").Value = 242.5
$ws.Range("N36").Value = -13430
$ws.Range("H70").Value = 1170.4375
$ws.Range("I70").Value = 746.7143
$ws.Range("J70").Value = 1500
$ws.Range("K70").Value = 2240.1429
$ws.Range("L70").Value = 4500
$ws.Range("M70").Value = -1970.1429
$ws.Range("N70").Value = -5040
$ws.Range("H73").Value = 1170.4375
$ws.Range("I73").Value = 746.7143
$ws.Range("J73").Value = 1500
$ws.Range("K73").Value = 2240.1429
$ws.Range("L73").Value = 4500
$ws.Range("M73").Value = -1304.1429
$ws.Range("N73").Value = -6372
$ws.Range("H74").Value = 2549972.5
$ws.Range("I74").Value = 3186534.5
$ws.Range("J74").Value = 3725
$ws.Range("K74").Value = 3186534.5
$ws.Range("L74").Value = 3725
$ws.Range("M74").Value = -3185598.5
$ws.Range("N74").Value = -5597
$ws.Range("H77").Value = 2549972.5
$ws.Range("I77").Value = 3186534.5
$ws.Range("J77").Value = 3725
$ws.Range("K77").Value = 15932672.5
$ws.Range("L77").Value = 18625
$ws.Range("M77").Value = -15927992.5
$ws.Range("N77").Value = -27985
$ws.Range("H80").Value = 278267.22
$ws.Range("I80").Value = 556.0625
$ws.Range("K80").Value = 1668.1875
$ws.Range("M80").Value = -670.1875
$ws.Range("H83").Value = 278267.22
$ws.Range("I83").Value = 556.0625
$ws.Range("K83").Value = 5004.5625
$ws.Range("M83").Value = -12.5625
$ws.Range("H86").Value = 55559450
$ws.Range("I86").Value = 2550.625
$ws.Range("J86").Value = 100004970
$ws.Range("K86").Value = 2550.625
$ws.Range("L86").Value = 100004970
$ws.Range("M86").Value = -1427.625
$ws.Range("N86").Value = -100007216
$ws.Range("H89").Value = 55559450
$ws.Range("I89").Value = 2550.625
$ws.Range("J89").Value = 100004970
$ws.Range("K89").Value = 12753.125
$ws.Range("L89").Value = 500024850
$ws.Range("M89").Value = -7137.125
$ws.Range("N89").Value = -500036082
$ws.Range("H103").Value = 843.0625
$ws.Range("I103").Value = 815.75
$ws.Range("J103").Value = 925
$ws.Range("K103").Value = 2447.25
$ws.Range("L103").Value = 2775
$ws.Range("M103").Value = -1861.25
$ws.Range("N103").Value = -3947
$ws.Range("H113").Value = 2660.5
$ws.Range("I113").Value = 2522.1428
$ws.Range("K113").Value = 2522.1428
$ws.Range("M113").Value = 731.8571999999999
$ws.Range("H137").Value = 888.5
$ws.Range("I137").Value = 832.75
$ws.Range("J137").Value = 1000
$ws.Range("K137").Value = 2498.25
$ws.Range("L137").Value = 3000
$ws.Range("M137").Value = 51.75
$ws.Range("N137").Value = -8100
$ws.Range("H138").Value = 2603.99
$ws.Range("I138").Value = 1073.1724
$ws.Range("J138").Value = 3229.2534
$ws.Range("K138").Value = 3219.5172
$ws.Range("L138").Value = 9687.760200000001
$ws.Range("M138").Value = 1920.4828
$ws.Range("N138").Value = -19967.7602

# Sheet: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Range("H74").Value = 1734.125
$ws.Range("I74").Value = 1631.4546
$ws.Range("J74").Value = 1960
$ws.Range("K74").Value = 1631.4546
$ws.Range("L74").Value = 1960
$ws.Range("M74").Value = -757.4546
$ws.Range("N74").Value = -3708
$ws.Range("H77").Value = 1734.125
$ws.Range("I77").Value = 1631.4546
$ws.Range("J77").Value = 1960
$ws.Range("K77").Value = 8157.273
$ws.Range("L77").Value = 9800
$ws.Range("M77").Value = -3789.273
$ws.Range("N77").Value = -18536
$ws.Range("H132").Value = 3227.8333
$ws.Range("I132").Value = 2822.4
$ws.Range("J132").Value = 3517.4285
$ws.Range("K132").Value = 8467.200000000001
$ws.Range("L132").Value = 10552.2855
$ws.Range("M132").Value = -5937.200000000001
$ws.Range("N132").Value = -15612.2855

# Sheet: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Range("H86").Value = 22729238
$ws.Range("I86").Value = 30304696
$ws.Range("K86").Value = 30304696
$ws.Range("M86").Value = -30303573
$ws.Range("H89").Value = 22729238
$ws.Range("I89").Value = 30304696
$ws.Range("K89").Value = 151523480
$ws.Range("M89").Value = -151517864
$ws.Range("H94").Value = 1761.4706
$ws.Range("I94").Value = 1754.6666
$ws.Range("J94").Value = 1777.8
$ws.Range("K94").Value = 1754.6666
$ws.Range("L94").Value = 1777.8
$ws.Range("M94").Value = -1303.6666
$ws.Range("N94").Value = -2679.8
$ws.Range("H99").Value = 30304798
$ws.Range("I99").Value = 50001740
$ws.Range("J99").Value = 1807.1538
$ws.Range("K99").Value = 50001740
$ws.Range("L99").Value = 1807.1538
$ws.Range("M99").Value = -50000242
$ws.Range("N99").Value = -4803.1538
$ws.Range("H105").Value = 5045
$ws.Range("I105").Value = 4239.2
$ws.Range("J105").Value = 7462.4
$ws.Range("K105").Value = 4239.2
$ws.Range("L105").Value = 7462.4
$ws.Range("M105").Value = -2492.2
$ws.Range("N105").Value = -10956.4

# Sheet: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 2596.7
$ws.Range("I31").Value = 2647.4412
$ws.Range("J31").Value = 2309.1667
$ws.Range("K31").Value = 2647.4412
$ws.Range("L31").Value = 2309.1667
$ws.Range("M31").Value = -2352.4412
$ws.Range("N31").Value = -2899.1667
$ws.Range("H34").Value = 2596.7
$ws.Range("I34").Value = 2647.4412
$ws.Range("J34").Value = 2309.1667
$ws.Range("K34").Value = 2647.4412
$ws.Range("L34").Value = 2309.1667
$ws.Range("M34").Value = -2445.4412
$ws.Range("N34").Value = -2713.1667
$ws.Range("H99").Value = 2613.0435
$ws.Range("I99").Value = 2600
$ws.Range("J99").Value = 2750
$ws.Range("K99").Value = 2600
$ws.Range("L99").Value = 2750
$ws.Range("M99").Value = -1102
$ws.Range("N99").Value = -5746
$ws.Range("H126").Value = 2613.0435
$ws.Range("I126").Value = 2600
$ws.Range("J126").Value = 2750
$ws.Range("K126").Value = 7800
$ws.Range("L126").Value = 8250
$ws.Range("M126").Value = -5330
$ws.Range("N126").Value = -13190
$ws.Range("H134").Value = 20834476
$ws.Range("I134").Value = 1091.0476
$ws.Range("K134").Value = 3273.142800000001
$ws.Range("M134").Value = -738.1428000000005

# Sheet: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents() | Out-Null

# Sheet: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Range("H126").Value = 4182.6
$ws.Range("I126").Value = 3821.5
$ws.Range("K126").Value = 11464.5
$ws.Range("M126").Value = -8994.5

# Sheet: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Range("H16").Value = 1443.1333
$ws.Range("I16").Value = 1559.4117
$ws.Range("J16").Value = 1291.0769
$ws.Range("K16").Value = 1559.4117
$ws.Range("L16").Value = 1291.0769
$ws.Range("M16").Value = -1389.4117
$ws.Range("N16").Value = -1631.0769
$ws.Range("H61").Value = 1672.091
$ws.Range("I61").Value = 1749.1765
$ws.Range("J61").Value = 1410
$ws.Range("K61").Value = 1749.1765
$ws.Range("L61").Value = 1410
$ws.Range("M61").Value = -1547.1765
$ws.Range("N61").Value = -1814
$ws.Range("H113").Value = 1672.091
$ws.Range("I113").Value = 1749.1765
$ws.Range("J113").Value = 1410
$ws.Range("K113").Value = 1749.1765
$ws.Range("L113").Value = 1410
$ws.Range("M113").Value = 420.8235
$ws.Range("N113").Value = -5750

# Sheet: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Range("H42").Value = 80049
$ws.Range("J42").Value = 80049
$ws.Range("L42").Value = 80049
$ws.Range("N42").Value = -80805
$ws.Range("H126").Value = 895.6667
$ws.Range("I126").Value = 895.6667
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 2687.0001
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -217.0001000000002
$ws.Range("N126").ClearContents() | Out-Null
$ws.Range("H132").Value = 1943.0625
$ws.Range("I132").Value = 1237.8462
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 3713.5386
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -1183.5386
$ws.Range("N132").Value = -20057
$ws.Range("H136").Value = 398.7586
$ws.Range("I136").Value = 420.5
$ws.Range("K136").Value = 1261.5
$ws.Range("M136").Value = 1288.5
